$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data (columns B through AC) between row pairs whose match
# records were re-ordered: (31,32), (188,189), (195,196). Column A (the
# row sequence id) is left untouched.
$pairs = @(@(31, 32), @(188, 189), @(195, 196))
foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $rangeA = $ws.Range("B" + $r1 + ":AC" + $r1)
    $rangeB = $ws.Range("B" + $r2 + ":AC" + $r2)
    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2
    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

# Remove the last four rows of data (rows 210-213), which are no longer
# present in the updated dataset.
$ws.Range("A210:A213").EntireRow.Delete()
